# Journal de travail - ajout d'une ligne d'activité au 10.02.2023
# (Analyse de la partie "gestion des tickets" sur la nouvelle version de
# l'application) et mise à jour des totaux en conséquence.

$wb = $excel.ActiveWorkbook

$wsJournal = $wb.Worksheets.Item("Journal")
$wsTotaux  = $wb.Worksheets.Item("Totaux")

# --- Journal: grow the table by two rows (1 data row + 1 trailing blank
# row, mirroring the existing layout) -------------------------------------
$tbl = $wsJournal.ListObjects.Item("Tableau1")
$tbl.ListRows.Add() | Out-Null
$tbl.ListRows.Add() | Out-Null

# Write the new entry's values first so dependent formulas recalc against
# the final data before the formatting-only paste below touches the cells.
$wsJournal.Range("A15").Value = 44967
$wsJournal.Range("B15").Value = 2
$wsJournal.Range("C15").Value = 0.0625
$wsJournal.Range("D15").Value = "Analyse"
$wsJournal.Range("E15").Value = "Analyse de la partie ""gestion des tickets"" sur la nouvelle version de l'application "

# Copy the formatting of the last pre-existing row onto the new rows so the
# cell styles (date/number formats, alignment) match the rest of the table.
$wsJournal.Range("A14:E14").Copy()
$wsJournal.Range("A15:E15").PasteSpecial(-4122)

$wsJournal.Range("A14:D14").Copy()
$wsJournal.Range("A16:D16").PasteSpecial(-4122)

# --- Totaux: extend the weekly sum to include the new Journal row --------
$wsTotaux.Range("B6").Formula = "=SUM(Journal!C13:C15)"

# --- Restore the selections recorded in the saved workbook ---------------
$wsJournal.Range("D29").Select()
$wsTotaux.Range("B10").Select()
